$wb = $excel.ActiveWorkbook

# --- Update sigma_010 (sheet 2) with refined values ---
$ws2 = $wb.Worksheets.Item("sigma_010")
$ws2.Range("B2").Value = 28.15749220272387
$ws2.Range("C2").Value = 33.25940235514033
$ws2.Range("B3").Value = 28.19896563328032
$ws2.Range("C3").Value = 33.25340765287542
$ws2.Range("B4").Value = 28.17809358717455
$ws2.Range("C4").Value = 33.26030823922135
$ws2.Range("B5").Value = 28.17980109862476
$ws2.Range("C5").Value = 33.26081629860633
$ws2.Range("B6").Value = 28.21414999486051
$ws2.Range("C6").Value = 33.25860301296565
$ws2.Range("B7").Value = 28.19725990262864
$ws2.Range("C7").Value = 33.24263964833596
$ws2.Range("B8").Value = 28.20993693052752
$ws2.Range("C8").Value = 33.24407200515023
$ws2.Range("B9").Value = 28.19908109582489
$ws2.Range("C9").Value = 33.26021059545074
$ws2.Range("B10").Value = 28.2194524876628
$ws2.Range("C10").Value = 33.25954096572669
$ws2.Range("B11").Value = 28.2189707392203
$ws2.Range("C11").Value = 33.24830439200367
$ws2.Range("B12").Value = 28.19732036725281
$ws2.Range("C12").Value = 33.25473051654764

# --- Update sigma_025 (sheet 3) with refined values ---
$ws3 = $wb.Worksheets.Item("sigma_025")
$ws3.Range("B2").Value = 19.85380687282416
$ws3.Range("C2").Value = 29.544222895223
$ws3.Range("B3").Value = 19.86470388157477
$ws3.Range("C3").Value = 29.52817161099905
$ws3.Range("B4").Value = 19.88884880910413
$ws3.Range("C4").Value = 29.59180374711842
$ws3.Range("B5").Value = 19.88856580198406
$ws3.Range("C5").Value = 29.51299154080481
$ws3.Range("B6").Value = 19.86271077357791
$ws3.Range("C6").Value = 29.59219087220933
$ws3.Range("B7").Value = 19.86807984186887
$ws3.Range("C7").Value = 29.56182093183921
$ws3.Range("B8").Value = 19.86832298211278
$ws3.Range("C8").Value = 29.54705751835244
$ws3.Range("B9").Value = 19.86591133274888
$ws3.Range("C9").Value = 29.46960142058028
$ws3.Range("B10").Value = 19.85172248307931
$ws3.Range("C10").Value = 29.53943119464385
$ws3.Range("B11").Value = 19.87156248773879
$ws3.Range("C11").Value = 29.54148089143867
$ws3.Range("B12").Value = 19.86842352666137
$ws3.Range("C12").Value = 29.54287726232091

# --- Add new sigma_050 sheet at the end with its results ---
$lastIndex = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($lastIndex)
$ws4 = $wb.Worksheets.Add($null, $afterSheet)
$ws4.Name = "sigma_050"

$ws4.Range("A1").Value = "Rows"
$ws4.Range("B1").Value = "Noisy"
$ws4.Range("C1").Value = "NLM-LBP"
$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = 14.9750551635481
$ws4.Range("C2").Value = 25.13139430844031
$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = 14.97644182520321
$ws4.Range("C3").Value = 25.08685774624304
$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = 14.98932466617925
$ws4.Range("C4").Value = 25.15067803761669
$ws4.Range("A5").Value = 3
$ws4.Range("B5").Value = 14.96603730031879
$ws4.Range("C5").Value = 25.10599678492797
$ws4.Range("A6").Value = 4
$ws4.Range("B6").Value = 14.99298740699105
$ws4.Range("C6").Value = 25.08144550619095
$ws4.Range("A7").Value = 5
$ws4.Range("B7").Value = 14.98230432369531
$ws4.Range("C7").Value = 25.09683889967036
$ws4.Range("A8").Value = 6
$ws4.Range("B8").Value = 14.98830347863831
$ws4.Range("C8").Value = 25.11233276670974
$ws4.Range("A9").Value = 7
$ws4.Range("B9").Value = 14.99666434806507
$ws4.Range("C9").Value = 25.1417399975558
$ws4.Range("A10").Value = 8
$ws4.Range("B10").Value = 15.01669892719094
$ws4.Range("C10").Value = 25.1613059939732
$ws4.Range("A11").Value = 9
$ws4.Range("B11").Value = 14.98581107539644
$ws4.Range("C11").Value = 25.11782533916428
$ws4.Range("A12").Value = "Média"
$ws4.Range("B12").Value = 14.98696285152265
$ws4.Range("C12").Value = 25.11864153804923
